# Update "想去人数" (number of people interested) counts on the
# "展览" (Exhibitions) and "全部类型" (All types) sheets to reflect the
# latest generated numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 160
$ws1.Range("F3").Value = 473
$ws1.Range("F4").Value = 13
$ws1.Range("F5").Value = 18
$ws1.Range("F8").Value = 15
$ws1.Range("F9").Value = 322

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 160
$ws4.Range("F4").Value = 473
$ws4.Range("F5").Value = 13
$ws4.Range("F6").Value = 18
$ws4.Range("F9").Value = 15
$ws4.Range("F10").Value = 322
